# Updates cryptos list prices/volumes (GitHub Actions style refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the cell to remain plain text even though the value looks numeric
    # (e.g. "1.000", "312.02"), mirroring the source data which stores these
    # as inline/shared strings rather than numbers.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# Row 2 - Bitcoin
Set-TextValue "D2" "27.710.83"
$ws.Range("E2").Value = "  -0.51%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.895.35"

# Row 4 - TetherUSD
Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  -1.16%  "

# Row 5 - BNB
Set-TextValue "D5" "312.02"
$ws.Range("E5").Value = "  -0.42%  "

# Row 6 - USDC
Set-TextValue "D6" "1.001"
$ws.Range("E6").Value = "  -1.02%  "

# Row 7 - XRP
Set-TextValue "D7" "0.4886"
$ws.Range("E7").Value = "  +1.18%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.3795"
$ws.Range("E8").Value = "  -0.55%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.07331"
$ws.Range("E9").Value = "  -0.48%  "

# Row 10 - Polygon
Set-TextValue "D10" "0.9129"
$ws.Range("E10").Value = "  -3.05%  "

# Row 11 - Solana
Set-TextValue "D11" "20.57"
$ws.Range("E11").Value = "  -2.18%  "

# Row 12 - TRON
Set-TextValue "D12" "0.07653"
$ws.Range("E12").Value = "  -1.69%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.915.47"
$ws.Range("E13").Value = "  +2.01%  "

# Row 14 - Polkadot
Set-TextValue "D14" "5.477"
$ws.Range("E14").Value = "  -0.73%  "

# Row 15 - Chainlink
Set-TextValue "D15" "6.618"
$ws.Range("E15").Value = "  +0.00%  "

# Row 16 - Litecoin
Set-TextValue "D16" "91.39"
$ws.Range("E16").Value = "  +0.05%  "

# Row 17 - BinanceUSD
$ws.Range("E17").Value = "  -1.18%  "

# Row 18 - ShibaInu
Set-TextValue "D18" "0.000008788"
$ws.Range("E18").Value = "  -0.80%  "

# Row 19 - Dai
Set-TextValue "D19" "1.000"
$ws.Range("E19").Value = "  -1.00%  "

# Row 20 - WrappedBTC
Set-TextValue "D20" "27.698.28"
$ws.Range("E20").Value = "  -0.65%  "

# Row 21 - Avalanche
Set-TextValue "D21" "14.50"
$ws.Range("E21").Value = "  -2.51%  "

# Row 22 - Uniswap
Set-TextValue "D22" "5.125"
$ws.Range("E22").Value = "  +0.00%  "

# Row 23 - WrappedliquidstakedEther2.0
Set-TextValue "D23" "2.186.41"
$ws.Range("E23").Value = "  +3.10%  "

# Row 24 - Cosmos
$ws.Range("E24").Value = "  -1.19%  "

# Row 25 - Toncoin
Set-TextValue "D25" "1.904"
$ws.Range("E25").Value = "  -2.27%  "

# Row 26 - Monero
Set-TextValue "D26" "153.98"
$ws.Range("E26").Value = "  -2.21%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "18.37"
$ws.Range("E27").Value = "  -1.07%  "

# Row 28 - LidoDAOToken
Set-TextValue "D28" "2.168"
$ws.Range("E28").Value = "  +5.98%  "

# Row 29 - BitcoinCash
Set-TextValue "D29" "115.41"
$ws.Range("E29").Value = "  -0.46%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextValue "D30" "4.889"
$ws.Range("E30").Value = "  -1.79%  "

# Row 31 - Stellar
Set-TextValue "D31" "0.08910"
$ws.Range("E31").Value = "  +0.28%  "

# Row 32 - HuobiToken
$ws.Range("E32").Value = "  -4.31%  "

# Row 33 - ARBITRUM
Set-TextValue "D33" "1.224"
$ws.Range("E33").Value = "  -0.03%  "

# Row 34 - ImmutableX
Set-TextValue "D34" "0.7682"
$ws.Range("E34").Value = "  -0.35%  "

# Row 35 - Filecoin
Set-TextValue "D35" "4.634"
$ws.Range("E35").Value = "  -0.44%  "

# Row 36 - now RenderToken (was VeChain)
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D36" "2.558"
$ws.Range("E36").Value = "  -6.65%  "

# Row 37 - now VeChain (was RenderToken)
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D37" "0.02037"
$ws.Range("E37").Value = "  -0.35%  "

# Row 38 - TrustWalletToken
$ws.Range("E38").Value = "  -2.71%  "

# Row 39 - Hedera
Set-TextValue "D39" "0.05278"
$ws.Range("E39").Value = "  -1.73%  "

# Row 40 - TheSandbox
Set-TextValue "D40" "0.5479"
$ws.Range("E40").Value = "  -2.66%  "

# Row 41 - MXToken
Set-TextValue "D41" "2.983"
$ws.Range("E41").Value = "  -0.56%  "

# Row 42 - FraxShare
Set-TextValue "D42" "6.892"
$ws.Range("E42").Value = "  -2.30%  "

# Row 43 - Aptos
Set-TextValue "D43" "8.547"
$ws.Range("E43").Value = "  -0.18%  "

# Row 44 - now Algorand (was Quant)
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D44" "0.1520"
$ws.Range("E44").Value = "  -0.70%  "

# Row 45 - now Quant (was Algorand)
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D45" "112.37"
$ws.Range("E45").Value = "  +7.08%  "

# Row 46 - EnergySwap
Set-TextValue "D46" "10.67"
$ws.Range("E46").Value = "  +0.27%  "

# Row 47 - Decentraland
Set-TextValue "D47" "0.4788"
$ws.Range("E47").Value = "  -1.92%  "

# Row 48 - PaxDollar
Set-TextValue "D48" "1.000"
$ws.Range("E48").Value = "  -1.12%  "

# Row 49 - NEARProtocol
Set-TextValue "D49" "1.641"
$ws.Range("E49").Value = "  -1.50%  "

# Row 50 - Aave
Set-TextValue "D50" "67.43"
$ws.Range("E50").Value = "  -0.97%  "

# Row 51 - Cronos
Set-TextValue "D51" "0.06049"
$ws.Range("E51").Value = "  -1.21%  "
